$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "100"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "200"
